# Scrape update: 30/12/2025 12:28 run, adding new rows scraped for
# LP1912 (sheet "LP1912"), and for 6203/6173 (sheet "6203-6173").
# "LP1912-215" only gets its "Última actualización" timestamp bumped
# (no new rows scraped this run).

$wb = $excel.ActiveWorkbook

$newTimestamp = "Última actualización: 30/12/2025 12:28:13"

# ---------------------------------------------------------------
# Sheet "LP1912": 18 new rows (204-221), totals 202 -> 220
# ---------------------------------------------------------------
$wsLP1912 = $wb.Worksheets.Item("LP1912")
$wsLP1912.Range("A2").Value = $newTimestamp
$wsLP1912.Range("A3").Value = "Total filas: 220"

$lp1912Rows = @(
    @("12:28:02", "12:35", "23_HERNANDEZ", 7, "LP1912", "30/12/2025"),
    @("12:28:02", "12:37", "16_SANTA ANA", 9, "LP1912", "30/12/2025"),
    @("12:28:02", "12:37", "27_EL RETIRO", 9, "LP1912", "30/12/2025"),
    @("12:28:02", "12:38", "17_179 Y 38", 10, "LP1912", "30/12/2025"),
    @("12:28:02", "12:47", "16_SANTA ANA", 19, "LP1912", "30/12/2025"),
    @("12:28:02", "12:50", "15_ABASTO", 22, "LP1912", "30/12/2025"),
    @("12:28:02", "12:55", "10_OLMOS", 27, "LP1912", "30/12/2025"),
    @("12:28:02", "13:02", "15_ABASTO", 34, "LP1912", "30/12/2025"),
    @("12:28:02", "13:07", "16_P MOR-SANTA ANA", 39, "LP1912", "30/12/2025"),
    @("12:28:02", "13:08", "10_OLMOS", 40, "LP1912", "30/12/2025"),
    @("12:28:02", "13:20", "10_OLMOS", 52, "LP1912", "30/12/2025"),
    @("12:28:02", "13:27", "14_ABASTO", 59, "LP1912", "30/12/2025"),
    @("12:28:02", "13:36", "15_ABASTO", 68, "LP1912", "30/12/2025"),
    @("12:28:02", "13:36", "23_HERNANDEZ", 68, "LP1912", "30/12/2025"),
    @("12:28:02", "13:46", "17_ROMERO", 78, "LP1912", "30/12/2025"),
    @("12:28:02", "13:57", "16_P MOR-167 Y 521", 89, "LP1912", "30/12/2025"),
    @("12:28:02", "14:04", "17_ROMERO", 96, "LP1912", "30/12/2025"),
    @("12:28:02", "14:05", "23_HERNANDEZ", 97, "LP1912", "30/12/2025")
)

$r = 204
foreach ($row in $lp1912Rows) {
    $wsLP1912.Cells.Item($r, 2).Value = $row[0]
    $wsLP1912.Cells.Item($r, 3).Value = $row[1]
    $wsLP1912.Cells.Item($r, 4).Value = $row[2]
    $wsLP1912.Cells.Item($r, 5).Value = $row[3]
    $wsLP1912.Cells.Item($r, 6).Value = $row[4]
    $wsLP1912.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "LP1912-215": no new rows this run, only timestamp bump
# ---------------------------------------------------------------
$wsLP1912215 = $wb.Worksheets.Item("LP1912-215")
$wsLP1912215.Range("A2").Value = $newTimestamp

# ---------------------------------------------------------------
# Sheet "6203-6173": 2 new rows (30-31), totals 28 -> 30
# ---------------------------------------------------------------
$ws6203 = $wb.Worksheets.Item("6203-6173")
$ws6203.Range("A2").Value = $newTimestamp
$ws6203.Range("A3").Value = "Total filas: 30"

$ws6203.Cells.Item(30, 2).Value = "30/12/2025"
$ws6203.Cells.Item(30, 3).Value = "12:28:08"
$ws6203.Cells.Item(30, 4).Value = "12:54"
$ws6203.Cells.Item(30, 5).Value = "215C_LA PLATA"
$ws6203.Cells.Item(30, 6).Value = 26
$ws6203.Cells.Item(30, 7).Value = "L6203"

$ws6203.Cells.Item(31, 2).Value = "30/12/2025"
$ws6203.Cells.Item(31, 3).Value = "12:28:13"
$ws6203.Cells.Item(31, 4).Value = "13:31"
$ws6203.Cells.Item(31, 5).Value = "215B_LP-P MOR-1 Y 57"
$ws6203.Cells.Item(31, 6).Value = 63
$ws6203.Cells.Item(31, 7).Value = "L6173"
